# Applies the per-cell odds/stake updates from the 2025-12-30 Betfair
# Back/Lay workbook diff. Sheet1 is a flat table (header row 1, data
# rows 2-33, columns A:AO) so each change is a direct cell-value write.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S2").Value = 1.43
$ws.Range("I3").Value = 1000
$ws.Range("K3").Value = 500
$ws.Range("N4").Value = 2.94
$ws.Range("Q4").Value = 1.83
$ws.Range("R4").Value = 1.25
$ws.Range("H6").Value = 1.26
$ws.Range("I6").Value = 1.41
$ws.Range("T6").Value = 1.61
$ws.Range("W6").Value = 1.05
$ws.Range("I7").Value = 4.3
$ws.Range("J7").Value = 3.4
$ws.Range("O7").Value = 1.37
$ws.Range("S7").Value = 2.18
$ws.Range("V7").Value = 1.33
$ws.Range("J8").Value = 5.9
$ws.Range("W8").Value = 4
$ws.Range("AC8").Value = 990
$ws.Range("AH8").Value = 990
$ws.Range("O9").Value = 1.43
$ws.Range("S9").Value = 4.4
$ws.Range("F10").Value = 5.6
$ws.Range("G10").Value = 5.9
$ws.Range("I10").Value = 1.69
$ws.Range("Q10").Value = 1.84
$ws.Range("R10").Value = 1.43
$ws.Range("T10").Value = 1.88
$ws.Range("U10").Value = 2.08
$ws.Range("V10").Value = 2.44
$ws.Range("W10").Value = 1.2
$ws.Range("X10").Value = 16.5
$ws.Range("AB10").Value = 20
$ws.Range("AN10").Value = 80
$ws.Range("AO10").Value = 9
$ws.Range("F11").Value = 1.6
$ws.Range("G11").Value = 1.61
$ws.Range("T11").Value = 1.69
$ws.Range("W11").Value = 2.64
$ws.Range("AD11").Value = 22
$ws.Range("G12").Value = 3.45
$ws.Range("H12").Value = 2.24
$ws.Range("I12").Value = 2.28
$ws.Range("L12").Value = 1.33
$ws.Range("P12").Value = 2.2
$ws.Range("V12").Value = 1.78
$ws.Range("AB12").Value = 15.5
$ws.Range("AF12").Value = 25
$ws.Range("AN12").Value = 29
$ws.Range("AO12").Value = 14.5
$ws.Range("N13").Value = 2.48
$ws.Range("P13").Value = 2.46
$ws.Range("Q13").Value = 1.44
$ws.Range("F14").Value = 2.24
$ws.Range("G14").Value = 3.2
$ws.Range("I14").Value = 3.85
$ws.Range("J14").Value = 2.8
$ws.Range("K14").Value = 6
$ws.Range("L14").Value = 1.35
$ws.Range("P14").Value = 1.62
$ws.Range("S14").Value = 2.66
$ws.Range("V14").Value = 1.35
$ws.Range("W14").Value = 1.46
$ws.Range("H15").Value = 2.78
$ws.Range("P15").Value = 1.9
$ws.Range("AB15").Value = 990
$ws.Range("AO15").Value = 34
$ws.Range("I16").Value = 7.2
$ws.Range("L16").Value = 1.37
$ws.Range("N16").Value = 2.6
$ws.Range("O16").Value = 1.29
$ws.Range("P16").Value = 1.75
$ws.Range("Q16").Value = 1.61
$ws.Range("U16").Value = 1.71
$ws.Range("Y16").Value = 990
$ws.Range("AH16").Value = 990
$ws.Range("AN16").Value = 13.5
$ws.Range("T17").Value = 1.76
$ws.Range("X17").Value = 13
$ws.Range("AN17").Value = 48
$ws.Range("AO17").Value = 38
$ws.Range("AN18").Value = 15
$ws.Range("F19").Value = 1.45
$ws.Range("N19").Value = 1.1
$ws.Range("T19").Value = 1.04
$ws.Range("U19").Value = 1.04
$ws.Range("W19").Value = 2.88
$ws.Range("AN19").Value = 6.2
$ws.Range("N20").Value = 2.44
$ws.Range("P20").Value = 2.06
$ws.Range("Q20").Value = 1.56
$ws.Range("R20").Value = 1.37
$ws.Range("S20").Value = 2.12
$ws.Range("Y20").Value = 990
$ws.Range("AD20").Value = 990
$ws.Range("AG20").Value = 990
$ws.Range("AH20").Value = 990
$ws.Range("AN20").Value = 9
$ws.Range("G21").Value = 1.41
$ws.Range("H21").Value = 9.2
$ws.Range("J21").Value = 5.3
$ws.Range("K21").Value = 15.5
$ws.Range("T21").Value = 1.04
$ws.Range("U21").Value = 1.04
$ws.Range("W21").Value = 3.35
$ws.Range("J22").Value = 3.65
$ws.Range("T22").Value = 1.6
$ws.Range("X22").Value = 990
$ws.Range("Y22").Value = 990
$ws.Range("AB22").Value = 990
$ws.Range("AC22").Value = 990
$ws.Range("AD22").Value = 990
$ws.Range("AG22").Value = 990
$ws.Range("AH22").Value = 990
$ws.Range("AN23").Value = 10.5
$ws.Range("I24").Value = 1.66
$ws.Range("F25").Value = 2.48
$ws.Range("N25").Value = 2.34
$ws.Range("P25").Value = 1.64
$ws.Range("Q25").Value = 1.96
$ws.Range("AD25").Value = 990
$ws.Range("AH25").Value = 990
$ws.Range("W26").Value = 2.6
$ws.Range("AN26").Value = 8.8
$ws.Range("H27").Value = 2.84
$ws.Range("AB27").Value = 17.5
$ws.Range("F28").Value = 1.9
$ws.Range("N28").Value = 2.62
$ws.Range("R28").Value = 1.22
$ws.Range("S28").Value = 2.52
$ws.Range("AN28").Value = 22
$ws.Range("T29").Value = 1.68
$ws.Range("AN29").Value = 29
$ws.Range("F30").Value = 3.1
$ws.Range("G30").Value = 5
$ws.Range("H30").Value = 1.8
$ws.Range("I30").Value = 1.98
$ws.Range("K30").Value = 980
$ws.Range("N30").Value = 2.36
$ws.Range("P30").Value = 2.2
$ws.Range("Q30").Value = 1.45
$ws.Range("R30").Value = 1.37
$ws.Range("S30").Value = 2.3
$ws.Range("T30").Value = 1.04
$ws.Range("U30").Value = 2.26
$ws.Range("W30").Value = 1.25
$ws.Range("F31").Value = 1.47
$ws.Range("G31").Value = 1.48
$ws.Range("M31").Value = 1.06
$ws.Range("R31").Value = 1.39
$ws.Range("X31").Value = 16
$ws.Range("Y31").Value = 24
$ws.Range("Z31").Value = 70
$ws.Range("AA31").Value = 300
$ws.Range("AH31").Value = 28
$ws.Range("AN31").Value = 7.6
$ws.Range("AO31").Value = 180
$ws.Range("Q32").Value = 1.63
$ws.Range("T32").Value = 1.95
$ws.Range("AA32").Value = 350
$ws.Range("G33").Value = 2.6
$ws.Range("H33").Value = 3.3
$ws.Range("L33").Value = 1.33
$ws.Range("N33").Value = 2.8
$ws.Range("O33").Value = 1.06
$ws.Range("Q33").Value = 1.95
$ws.Range("T33").Value = 1.04
$ws.Range("U33").Value = 1.04
$ws.Range("J31").Value = 4.9
$ws.Range("K31").Value = 5
